# Update NATMI LR-pair output (Spp1-Cd44) with newly recomputed TPM-based
# statistics. The ligand-side stats (G:J) depend only on the sending
# cluster, the receptor-side stats (M:P) depend only on the target
# cluster, and the edge-derived stats (Q:T) are specific to each
# sending/target pair (i.e. each data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligand-side values (columns G,H,I,J), keyed by "Sending cluster" (col A)
$ligandBySender = @{
    "ECs"           = @(0.7476426666666667, 2.242928, 0.001581772089386036, 0.001581772089386036)
    "FAPs"          = @(12.24662533333333, 36.739876, 0.02590993131491687, 0.02590993131491688)
    "MuSCs"         = @(18.93023433333333, 56.79070299999999, 0.04005030430848061, 0.04005030430848062)
    "Resolving-Mac" = @(440.7369333333333, 1322.2108, 0.9324579922872165, 0.9324579922872166)
}

# Receptor-side values (columns M,N,O,P), keyed by "Target cluster" (col D)
$receptorByTarget = @{
    "ECs"           = @(24.576554, 73.729662, 0.07553767049546639, 0.07553767049546638)
    "FAPs"          = @(117.0512696666667, 351.153809, 0.359764849016532, 0.359764849016532)
    "MuSCs"         = @(55.68784966666667, 167.063549, 0.1711603033819035, 0.1711603033819035)
    "Resolving-Mac" = @(128.0392633333333, 384.11779, 0.3935371771060981, 0.3935371771060981)
}

# Edge-derived values (columns Q,R,S,T), specific to each data row (2..17)
$edgeByRow = @{
    2  = @(18.37448037003734, 165.370323330336, 0.0001194833788869678, 0.0001194833788869678)
    3  = @(87.51252339030579, 787.6127105127521, 0.0005690659969165315, 0.0005690659969165316)
    4  = @(41.63461242571912, 374.7115118314721, 0.0002707365907003413, 0.0002707365907003413)
    5  = @(95.7276162765689, 861.5485464891201, 0.0006224861228821952, 0.0006224861228821953)
    6  = @(300.9798488224347, 2708.818639401912, 0.001957175854226357, 0.001957175854226357)
    7  = @(1433.483044398632, 12901.34739958769, 0.009321482527539783, 0.009321482527539785)
    8  = @(681.9882304866584, 6137.894074379925, 0.004434751704465455, 0.004434751704465455)
    9  = @(1568.048885999338, 14112.43997399404, 0.01019652122868528, 0.01019652122868528)
    10 = @(465.2399263258206, 4187.159336932385, 0.003025306690097166, 0.003025306690097166)
    11 = @(2215.807963804192, 19942.27167423773, 0.01440869168260669, 0.01440869168260669)
    12 = @(1054.184043709439, 9487.656393384947, 0.0068550222359771, 0.0068550222359771)
    13 = @(2423.813258767374, 21814.31932890637, 0.01576128369979966, 0.01576128369979966)
    14 = @(10831.79504186107, 97486.15537674959, 0.0704357045722559, 0.0704357045722559)
    15 = @(51588.81763565969, 464299.3587209372, 0.335465608809469, 0.335465608809469)
    16 = @(24543.69208601435, 220893.2287741292, 0.1595997928507606, 0.1595997928507606)
    17 = @(56431.63226779244, 507884.6904101319, 0.366956886054731, 0.366956886054731)
}

# Column numbers: G=7 H=8 I=9 J=10 M=13 N=14 O=15 P=16 Q=17 R=18 S=19 T=20
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $sender = $ws.Cells.Item($r, 1).Value()
    $target = $ws.Cells.Item($r, 4).Value()

    $lig = $ligandBySender[$sender]
    if ($lig) {
        $ws.Cells.Item($r, 7).Value  = $lig[0]
        $ws.Cells.Item($r, 8).Value  = $lig[1]
        $ws.Cells.Item($r, 9).Value  = $lig[2]
        $ws.Cells.Item($r, 10).Value = $lig[3]
    }

    $rec = $receptorByTarget[$target]
    if ($rec) {
        $ws.Cells.Item($r, 13).Value = $rec[0]
        $ws.Cells.Item($r, 14).Value = $rec[1]
        $ws.Cells.Item($r, 15).Value = $rec[2]
        $ws.Cells.Item($r, 16).Value = $rec[3]
    }

    $edge = $edgeByRow[$r]
    if ($edge) {
        $ws.Cells.Item($r, 17).Value = $edge[0]
        $ws.Cells.Item($r, 18).Value = $edge[1]
        $ws.Cells.Item($r, 19).Value = $edge[2]
        $ws.Cells.Item($r, 20).Value = $edge[3]
    }
}
